# Update the "Förändrad" date column (C) for rows 2 through 23
# from 2023-09-06 (serial 45175) to 2023-09-14 (serial 45183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45183
    }
}
